# Apply the edits described by the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet references ---
$wsData = $wb.Worksheets.Item("Data Integrity Testing ")
$wsCss  = $wb.Worksheets.Item("CSS-Accept-DB")

# --- 1. Rename the "Data Integrity Testing " sheet to "Data-Integrity-Testing" ---
$wsData.Name = "Data-Integrity-Testing"

# --- 2. Data-Integrity-Testing sheet: column widths ---
# (Input values are chosen so the engine's internal pixel-quantized column
# width lands on the closest representable value to the target OOXML
# `width` attribute - the ColumnWidth COM property and the raw stored
# width are not in 1:1 units in this engine.)
$wsData.Columns.Item(1).ColumnWidth = 12.333333333333332
$wsData.Columns.Item(2).ColumnWidth = 4.0
$wsData.Columns.Item(3).ColumnWidth = 9.166666666666668
$wsData.Columns.Item(4).ColumnWidth = 19.666666666666664
$wsData.Columns.Item(5).ColumnWidth = 22.166666666666664
$wsData.Columns.Item(6).ColumnWidth = 64.33333333333334
$wsData.Columns.Item(7).ColumnWidth = 58.83333333333333

# --- 3. Row heights ---
$wsData.Rows.Item(2).RowHeight = 30
$wsData.Rows.Item(3).RowHeight = 45
$wsData.Rows.Item(4).RowHeight = 45

# --- 4. Wrap text for F2:F4 (style s=14 -> s=15) ---
$wsData.Range("F2:F4").WrapText = $true

# --- 5. D4 new value "Verify record" (new shared string at index 41) ---
$wsData.Range("D4").Value = "Verify record"

# --- 6. View settings for Data-Integrity-Testing sheet ---
$wsData.Application.ActiveWindow.Zoom = 130
$wsData.Range("E11").Select()

# --- 7. View settings for CSS-Accept-DB sheet ---
$wsCss.Activate()
$wsCss.Range("E25").Select()

# Reactivate the renamed data sheet (it was tabSelected originally)
$wsData.Activate()
